# "added in Soft Duplicate check and KI database"
#
# The field team work-distribution sheet is re-synced against an updated
# contact/KI (Key Informant) roster:
#   - the "Responsible FO" column header is renamed to "Responsible_FO"
#     (snake_cased, presumably for a downstream script/field-mapping)
#   - several "Responsible FO" names are corrected / reassigned to the
#     current focal points
#   - two Banadir location labels are disambiguated with a "Mogadishu"
#     prefix ("Dayniile" -> "Mogadishu Dayniile", "Khada" -> "Mogadishu
#     Khada") to avoid collisions with identically named KI database
#     entries elsewhere ("Soft Duplicate check")
#   - row 66 ("Xudur") is moved from the Somaliland/"Omar Abdikarin" team
#     to the Southwest/"Abdikani/Hassan" team

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header rename -----------------------------------------------------
$ws.Range("D1").Value = "Responsible_FO"

# --- Banadir location disambiguation -----------------------------------
$ws.Range("C28").Value = "Mogadishu Dayniile"
$ws.Range("C46").Value = "Mogadishu Khada"

# --- row 66 ("Xudur") reassigned from Somaliland to Southwest ----------
$ws.Range("B66").Value = "Southwest"

# --- Responsible FO (column D) focal-point updates ---------------------
$ws.Range("D2").Value  = "Kunow/Hassan"
$ws.Range("D3").Value  = "Mohamed Kala"
$ws.Range("D4").Value  = "Mohamed Kala"
$ws.Range("D5").Value  = "Kunow/Hassan"
$ws.Range("D6").Value  = "Mohamed Hassan"
$ws.Range("D8").Value  = "Kunow/Hassan"
$ws.Range("D10").Value = "Mohamed Kala"
$ws.Range("D11").Value = "Mohamed Hassan"
$ws.Range("D12").Value = "Suleiman Mohamed"
$ws.Range("D13").Value = "Mohamed Hassan"
$ws.Range("D18").Value = "Mohamed Hassan"
$ws.Range("D19").Value = "Kunow/Hassan"
$ws.Range("D20").Value = "Suleiman Mohamed"
$ws.Range("D21").Value = "Isse Ahad"
$ws.Range("D22").Value = "Isse Ahad"
$ws.Range("D25").Value = "Kunow/Hassan"
$ws.Range("D26").Value = "Mohamed Kala"
$ws.Range("D28").Value = "Suleiman Mohamed"
$ws.Range("D29").Value = "Suleiman Mohamed"
$ws.Range("D30").Value = "Isse Ahad"
$ws.Range("D31").Value = "Kunow/Hassan"
$ws.Range("D32").Value = "Mohamed Kala"
$ws.Range("D33").Value = "Omar"
$ws.Range("D35").Value = "Mohamed Kala"
$ws.Range("D37").Value = "Suleiman Mohamed"
$ws.Range("D38").Value = "Suleiman Mohamed"
$ws.Range("D39").Value = "Mohamed Hassan"
$ws.Range("D40").Value = "Suleiman Mohamed"
$ws.Range("D41").Value = "Omar"
$ws.Range("D42").Value = "Suleiman Mohamed"
$ws.Range("D45").Value = "Suleiman Mohamed"
$ws.Range("D46").Value = "Suleiman Mohamed"
$ws.Range("D47").Value = "Mohamed Kala"
$ws.Range("D49").Value = "Mohamed Hassan"
$ws.Range("D50").Value = "Mohamed Hassan"
$ws.Range("D51").Value = "Mohamed Kala"
$ws.Range("D52").Value = "Kunow/Hassan"
$ws.Range("D54").Value = "Kunow/Hassan"
$ws.Range("D56").Value = "Suleiman Mohamed"
$ws.Range("D58").Value = "Suleiman Mohamed"
$ws.Range("D59").Value = "Mohamed Hassan"
$ws.Range("D60").Value = "Suleiman Mohamed"
$ws.Range("D61").Value = "Kunow/Hassan"
$ws.Range("D62").Value = "Suleiman Mohamed"
$ws.Range("D63").Value = "Kunow/Hassan"
$ws.Range("D64").Value = "Suleiman Mohamed"
$ws.Range("D65").Value = "Mohamed Hassan"
$ws.Range("D66").Value = "Abdikani/Hassan"
$ws.Range("D67").Value = "Suleiman Mohamed"
$ws.Range("D68").Value = "Mohamed Hassan"

# --- print setup (page orientation explicitly set to Portrait) ---------
$ws.PageSetup.Orientation = 1

# --- window/view state: scroll to row 34 and select C37 ----------------
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C37").Select()
